# Commit: "add account information when create a new role add init property module"
#
# Semantic edits to the Property1 sheet of EffectData.xlsx:
#  1. The five "gate" columns (DIZZY_GATE / MOVE_GATE / SKILL_GATE /
#     PHYSICAL_GATE / MAGIC_GATE -> AB:AF) are reset from 50 to 0 for every
#     data row (11-70).
#  2. Column AG ("BUFF_GATE") gets its own explicit width (it used to share
#     a merged width definition with AF) - the user widened/best-fit it.
#  3. The frozen-pane view scrolled right and the active selection moved to
#     the AB column (the first of the columns that was just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Reset DIZZY_GATE..MAGIC_GATE (AB:AF) to 0 for rows 11-70 in one shot.
$ws.Range("AB11:AF70").Value = 0

# 2) Give column AG its own width (previously AF:AG shared one column
#    definition); closest reachable width to the authored 15.73046875.
$ws.Columns("AG").ColumnWidth = 15

# 3) Move the selection onto the newly-edited column so AB11:AB70 is
#    highlighted, matching what a user would select after editing it.
$ws.Range("AB11:AB70").Select() | Out-Null
